$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.446.43'
$ws.Range('E2').Value = '  +5.22%  '
$ws.Range('D3').Value = '1.722.56'
$ws.Range('E3').Value = '  +4.54%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').Value = '  +0.26%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '225.50'
$ws.Range('E5').Value = '  +3.31%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5342'
$ws.Range('E6').Value = '  +2.77%  '
$ws.Range('E7').Value = '  +0.15%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2657'
$ws.Range('E8').Value = '  +1.56%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06567'
$ws.Range('E9').Value = '  +4.32%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.53'
$ws.Range('E10').Value = '  +6.05%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07669'
$ws.Range('E11').Value = '  +0.25%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.584'
$ws.Range('E12').Value = '  +0.29%  '
$ws.Range('D13').Value = '1.726.64'
$ws.Range('E13').Value = '  +3.88%  '
$ws.Range('D14').Value = '1.963.01'
$ws.Range('E14').Value = '  +4.78%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5783'
$ws.Range('E15').Value = '  +3.80%  '
$ws.Range('D16').Value = '0.0₅8273'
$ws.Range('E16').Value = '  +2.10%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '67.76'
$ws.Range('E17').Value = '  +4.28%  '
$ws.Range('D18').Value = '27.497.02'
$ws.Range('E18').Value = '  +5.63%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '218.14'
$ws.Range('E19').Value = '  +12.76%  '
$ws.Range('E20').Value = '  +0.09%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.717'
$ws.Range('E21').Value = '  +2.70%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.56'
$ws.Range('E22').Value = '  +1.34%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.031'
$ws.Range('E23').Value = '  +2.01%  '
$ws.Range('E24').Value = '  +0.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '142.85'
$ws.Range('E25').Value = '  -1.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.750'
$ws.Range('E26').Value = '  +16.54%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1231'
$ws.Range('E27').Value = '  +4.48%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.328'
$ws.Range('E28').Value = '  +2.04%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '16.46'
$ws.Range('E29').Value = '  +4.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05484'
$ws.Range('E30').Value = '  +1.68%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.301'
$ws.Range('E31').Value = '  +2.65%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.557'
$ws.Range('E32').Value = '  +3.42%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.433'
$ws.Range('E33').Value = '  +3.37%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.655'
$ws.Range('E34').Value = '  +6.51%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.855'
$ws.Range('E35').Value = '  +2.73%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9565'
$ws.Range('E36').Value = '  +1.71%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.425'
$ws.Range('E37').Value = '  +0.44%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5936'
$ws.Range('E38').Value = '  +6.56%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01646'
$ws.Range('E39').Value = '  +4.71%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.899'
$ws.Range('E40').Value = '  +2.64%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8493'
$ws.Range('E41').Value = '  +3.06%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '1.047.34'
$ws.Range('E42').Value = '  +1.99%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.004'
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.15'
$ws.Range('E44').Value = '  +0.46%  '
$ws.Range('D45').Value = '1.867.17'
$ws.Range('E45').Value = '  +4.72%  '
$ws.Range('E46').Value = '  +5.52%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '58.63'
$ws.Range('E47').Value = '  +2.52%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4482'
$ws.Range('E48').Value = '  +3.89%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.185'
$ws.Range('E49').Value = '  +3.93%  '
$ws.Range('E50').Value = '  +0.34%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05248'
$ws.Range('E51').Value = '  +3.01%  '
